# Actualización automática 2025-07-03 12:35:11
# Updates the "CUMPLIMIENTO MENSUAL" sheet: resets actual sales (VENTA) to 0
# for most product groups (so POR CUMPLIR == PRESUPUESTO and CUMPLIMIENTO == 0),
# applies a few PRESUPUESTO corrections, keeps a reduced VENTA figure for the
# PORCELANATO group, recomputes the TOTAL row, and tweaks a few column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Column width adjustments (D, E, F) -------------------------------------
# This runtime stores OOXML column <col width> as ColumnWidth + 0.8333333333333334,
# so we subtract that offset to land on the exact target widths (13, 22, 25).
$colOffset = 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 13 - $colOffset
$ws.Columns.Item(5).ColumnWidth = 22 - $colOffset
$ws.Columns.Item(6).ColumnWidth = 25 - $colOffset

# --- Row 2: 240X120 PORCELANATO ---------------------------------------------
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 9970.34304517915
$ws.Range("F2").Value = 0

# --- Row 3: 240X80 PORCELANATO ----------------------------------------------
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 27457.0076
$ws.Range("F3").Value = 0

# --- Row 4: FREGADEROS DE COCINA --------------------------------------------
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1003
$ws.Range("F4").Value = 0

# --- Row 5: GRANITO -----------------------------------------------------------
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 238.32
$ws.Range("F5").Value = 0

# --- Row 6: GRIFERIAS ---------------------------------------------------------
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 106.82
$ws.Range("F6").Value = 0

# --- Row 7: INODOROS (PRESUPUESTO also changes 1400 -> 2400) -----------------
$ws.Range("C7").Value = 2400
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 2400
$ws.Range("F7").Value = 0

# --- Row 8: LAVABOS ------------------------------------------------------------
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1000
$ws.Range("F8").Value = 0

# --- Row 9: LED -----------------------------------------------------------------
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 300
$ws.Range("F9").Value = 0

# --- Row 10: NO RESURTIBLES ------------------------------------------------------
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1300.5
$ws.Range("F10").Value = 0

# --- Row 11: OTROS (unchanged) ---------------------------------------------------

# --- Row 12: PANELES DECORATIVOS --------------------------------------------------
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 350
$ws.Range("F12").Value = 0

# --- Row 13: PANELES PU (PRESUPUESTO 230 -> 130) ----------------------------------
$ws.Range("C13").Value = 130
$ws.Range("E13").Value = 130

# --- Row 14: PANELES PVC (PRESUPUESTO 966 -> 240) ---------------------------------
$ws.Range("C14").Value = 240
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 240
$ws.Range("F14").Value = 0

# --- Row 15: PIEDRA SINTERIZADA ----------------------------------------------------
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 13500
$ws.Range("F15").Value = 0

# --- Row 16: PORCELANATO (PRESUPUESTO 32741.45 -> 51826.46, VENTA 30753.54 -> 4482.48) --
$ws.Range("C16").Value = 51826.46
$ws.Range("D16").Value = 4482.48
$ws.Range("E16").Value = 47343.98
$ws.Range("F16").Value = 0.08649018281395256

# --- Row 17: PUERTAS DE SEGURIDAD (unchanged) ---------------------------------------

# --- Row 18: SAL SOLUBLE --------------------------------------------------------------
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 3200
$ws.Range("F18").Value = 0

# --- Row 19: TOTAL (recomputed sums) --------------------------------------------------
$ws.Range("C19").Value = 113706.4506451791
$ws.Range("D19").Value = 4482.48
$ws.Range("E19").Value = 109223.9706451792
$ws.Range("F19").Value = 0.03942151016557164
